$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new module rows ---------------------------------------
# Row 11: "Android App(Basic) module" -- pushes the existing
#   "Notification module" row (and TOTAL / ratio rows below it) down by one.
$ws.Rows("11:11").Insert()

# Row 13 (after the first shift, this sits right below the now-shifted
#   "Notification module" row and above "TOTAL") : "Android Notification
#   Display Module" -- pushes TOTAL / ratio rows down by one more.
$ws.Rows("13:13").Insert()

# --- Populate the new rows ---------------------------------------------
$ws.Range("A11").Value = "Android App(Basic) module"
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0

$ws.Range("A13").Value = "Android Notification Display Module"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# --- Fix up the TOTAL row (now row 14) formulas so they cover the
#     extended data range (rows 3 through 13) -----------------------------
$ws.Range("B14").Formula = "=SUM(B3:B13)"
$ws.Range("C14").Formula = "=SUM(C3:C13)"
$ws.Range("D14").Formula = "=SUM(D3:D13)"
$ws.Range("E14").Formula = "=SUM(E3:E13)"
$ws.Range("F14").Formula = "=SUM(F3:F13)"
$ws.Range("G14").Formula = "=SUM(G3:G13)"

# --- Defect-ratio rows (now rows 16 and 17) keep referencing the TOTAL
#     row, which has simply moved from row 12 to row 14 ------------------
$ws.Range("C16").Formula = "=D14/B14"
$ws.Range("C17").Formula = "=G14/E14"

# --- Update the active selection to match the edited workbook -----------
$ws.Range("E18").Select()

Write-Output "done"
